# Retrieve price / product rows and write them into the sheet, checking
# whether the destination cell is already occupied before writing to it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old scratch data ("Hello World" / "Badland" / "Test1") is no longer
# needed now that column A is used to report occupancy, so clear it out.
$ws.Cells.Item(1, 6).Value = ""
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(3, 6).Value = ""
$ws.Cells.Item(6, 3).Value = ""

# Reserve the first 7 rows of column A for the occupied-slot notice.
for ($row = 1; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $existing = $cell.Value()
    if (-not [string]::IsNullOrEmpty($existing)) {
        # Slot already had data in it -> flag it as occupied.
        $cell.Value = "This is Occupied"
    } else {
        # Slot was free -> still mark it occupied now that we claim it.
        $cell.Value = "This is Occupied"
    }
}

# Walk one row past the reserved block: it is free, so report the miss.
$lastCell = $ws.Cells.Item(8, 1)
$lastExisting = $lastCell.Value()
if ([string]::IsNullOrEmpty($lastExisting)) {
    $lastCell.Value = "You have been juke"
} else {
    $lastCell.Value = "This is Occupied"
}

# Park the selection just below the rows we touched.
$ws.Range("A9:A12").Select() | Out-Null
